$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the formatting of the previous entry (row 26) down into the new
# row 27 - this is how a new leetcode entry was added to the tracker.
$ws.Range("A26:D26").Copy() | Out-Null
$ws.Range("A27:D27").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A27").Value = "Binary Search"
$ws.Range("B27").Value = "33. Search in Rotated Sorted Array"
$ws.Range("C27").Value = "for this question u need to know problem - Minimum in Rotated Sorted Array - Binary Search - Leetcode 153`nthe minimum in rotated array is nothing but the border/pivot element in rotated arr, so once u get its indx u can do simple binary search on the 2 subarray arrays divided by pivot"

# Unlike row 26 (which has a hyperlink in column D), row 27 has no entry in
# column D - remove the formatting/content that PasteSpecial carried over.
$ws.Range("D27").Clear() | Out-Null

# Row height follows from the wrapped 3-line summary, matching the other
# 3-line rows (5, 7, 9, 10, 16, 19) in the sheet.
$ws.Rows.Item(27).RowHeight = 43.2

# Scroll the view down to the newly added row and select the title cell,
# mirroring where the editor ended up after typing the new entry.
$ws.Application.Goto($ws.Range("A26"))
$ws.Range("B27").Select() | Out-Null
